$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.763.20'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').Value = '3.500.34'
$ws.Range('E3').Value = '  -3.63%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '578.48'
$ws.Range('E5').Value = '  -4.61%  '
$ws.Range('D6').Value = '192.49'
$ws.Range('E6').Value = '  -4.02%  '
$ws.Range('E7').Value = '  -2.30%  '
$ws.Range('D8').Value = '3.489.01'
$ws.Range('E8').Value = '  -3.62%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  -5.63%  '
$ws.Range('D11').Value = '0.622'
$ws.Range('E11').Value = '  -4.05%  '
$ws.Range('D12').Value = '51.53'
$ws.Range('E12').Value = '  -4.48%  '
$ws.Range('E13').Value = '  -6.30%  '
$ws.Range('D14').Value = '9.17'
$ws.Range('D15').Value = '4.055.93'
$ws.Range('E15').Value = '  -3.64%  '
$ws.Range('D16').Value = '644.91'
$ws.Range('E16').Value = '  -5.06%  '
$ws.Range('D17').Value = '69.680.77'
$ws.Range('E17').Value = '  -1.80%  '
$ws.Range('D18').Value = '3.494.10'
$ws.Range('E18').Value = '  -3.43%  '
$ws.Range('D19').Value = '12.44'
$ws.Range('E19').Value = '  -3.76%  '
$ws.Range('E20').Value = '  -1.81%  '
$ws.Range('D21').Value = '18.37'
$ws.Range('E21').Value = '  -3.48%  '
$ws.Range('D22').Value = '0.951'
$ws.Range('E22').Value = '  -4.94%  '
$ws.Range('D23').Value = '18.06'
$ws.Range('E23').Value = '  -2.50%  '
$ws.Range('D24').Value = '5.35'
$ws.Range('E24').Value = '  -0.78%  '
$ws.Range('D25').Value = '99.03'
$ws.Range('E25').Value = '  -5.86%  '
$ws.Range('E26').Value = '  -7.37%  '
$ws.Range('E27').Value = '  -4.44%  '
$ws.Range('D28').Value = '10.09'
$ws.Range('E28').Value = '  -4.07%  '
$ws.Range('D29').Value = '9.41'
$ws.Range('E29').Value = '  -4.32%  '
$ws.Range('D30').Value = '32.67'
$ws.Range('E30').Value = '  -4.59%  '
$ws.Range('D31').Value = '4.24'
$ws.Range('E31').Value = '  -8.32%  '
$ws.Range('D32').Value = '6.74'
$ws.Range('E32').Value = '  -6.31%  '
$ws.Range('E33').Value = '  -4.32%  '
$ws.Range('E34').Value = '  -4.92%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '61.50'
$ws.Range('E35').Value = '  -2.87%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').Value = '563.75'
$ws.Range('E36').Value = '  +9.26%  '
$ws.Range('B37').Value = 'CoreDAO'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D37').Value = '3.92'
$ws.Range('E37').Value = '  +51.73%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '3.704.07'
$ws.Range('E38').Value = '  -6.61%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0792'
$ws.Range('E40').Value = '  -9.05%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '3.62'
$ws.Range('E41').Value = '  +2.08%  '
$ws.Range('D42').Value = '2.91'
$ws.Range('E42').Value = '  -3.48%  '
$ws.Range('E43').Value = '  -3.61%  '
$ws.Range('E44').Value = '  -1.66%  '
$ws.Range('D45').Value = '34.34'
$ws.Range('E45').Value = '  -6.05%  '
$ws.Range('E46').Value = '  -3.45%  '
$ws.Range('E47').Value = '  -3.24%  '
$ws.Range('E48').Value = '  -6.98%  '
$ws.Range('E49').Value = '  -4.27%  '
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('D51').Value = '8.20'
$ws.Range('E51').Value = '  -5.21%  '
